# Remove the "W" oval (the compressor marker) from the DHP diagram slide --
# the sensor described by that callout isn't actually mounted at the
# compressor, so the shape needs to go.
#
# (The deck's datetimeFigureOut placeholders on the layouts/master also
# drift between 11/21/2013 and 12/5/2013 in the source history, but that is
# just PowerPoint re-caching the "today" field on whatever date the file was
# last saved -- not a deliberate content edit, and re-typing it here would
# only replace the live <a:fld> with static text. So we leave those alone
# and only perform the actual authored change: deleting the shape.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Oval 68") {
        $target = $sh
        break
    }
}

if ($target -eq $null) {
    # Fallback: locate by its cached shape id (69) or by its text "W" inside
    # a small ellipse, in case the name ever differs.
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Id -eq 69) {
            $target = $sh
            break
        }
    }
}

if ($target -ne $null) {
    $target.Delete()
}
